$d = $word.ActiveDocument

# --- Section 1: add "Second line of first section." and an empty paragraph
#     right after the existing "First section text." paragraph, without
#     disturbing that paragraph's own runs. ---
$p1 = $d.Paragraphs(1)
$ip = $p1.Range.Duplicate()
$ip.Collapse(0)               # wdCollapseEnd
$ip.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "Second line of first section."

# Split off a trailing empty paragraph (same style) using Find/Replace with a
# paragraph mark so no stray empty run is left behind in the new paragraph.
$d.Paragraphs(2).Range.Find.Execute("Second line of first section.", $true, $false, $false, $false, $false, $true, 1, $false, "Second line of first section.^p", 2)

# --- Section 2: turn the "Second section on next page." paragraph into the
#     end-of-section marker for a new (3rd) section by giving it its own
#     sectPr (continuous break), then add the new last-section paragraph. ---
$target = $d.Paragraphs($d.Paragraphs.Count)
$ip2 = $target.Range.Duplicate()
$ip2.Collapse(0)              # wdCollapseEnd
$ip2.InsertBreak(3)           # wdSectionBreakContinuous

# The break above created a fresh trailing paragraph carrying the new
# sectPr; merge it back into the target paragraph so the sectPr lives on
# that same paragraph (matching the target structure) instead of a
# separate empty one.
$target = $d.Paragraphs($d.Paragraphs.Count - 1)
$mark = $d.Range($target.Range.End - 1, $target.Range.End)
$mark.Delete()

# Add the final paragraph belonging to the new (continuous) last section.
$target = $d.Paragraphs($d.Paragraphs.Count)
$ip3 = $target.Range.Duplicate()
$ip3.Collapse(0)              # wdCollapseEnd
$ip3.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = "Last section on the same page."

Write-Output "done"
